$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Update data rows (row 2 through row 11)
$data = @(
    @(0, 3402.533333333333, 3671, 3175, 0.05813307762145996),
    @(1, 3593.833333333333, 3855, 3257, 0.05616005261739095),
    @(2, 3688.033333333333, 3963, 3381, 0.06020108064015706),
    @(3, 3503.7,            3755, 3143, 0.05842487017313639),
    @(4, 2647.833333333333, 2933, 2307, 0.06074936389923095),
    @(5, 2830.033333333333, 3060, 2460, 0.0578916072845459),
    @(6, 3603.9,            3888, 3214, 0.06155671278635661),
    @(7, 3156.3,            3438, 2789, 0.06042404174804687),
    @(8, 3470.733333333333, 3720, 3117, 0.0595992644627889),
    @(9, 3068.066666666667, 3354, 2789, 0.05703778266906738)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
